$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert the new "e036 - Re-enter Friendly Control" event as row 37 (pushing
# the existing rows 37-41 down to 38-42, matching e036's place right after
# e035 and before the end-game rows).
# ---------------------------------------------------------------------------
$ws.Rows.Item(37).Insert()

$ws.Range("A37").Value = "e036"

$bodyText = @'
<Bold>e036 Re-enter Friendly Control</Bold> 
<InlineUIContainer><Button Content='r4.54.5' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
No combat. 
<LineBreak/><LineBreak/>
'@
$ws.Range("B37").Value = $bodyText

# ---------------------------------------------------------------------------
# Row heights were recalculated (Excel re-wrapped the text) across the whole
# sheet after the edit. Re-apply the new heights for every data row.
# ---------------------------------------------------------------------------
$rowHeights = @{
    1  = 158.4
    2  = 187.2
    3  = 115.2
    4  = 100.8
    5  = 100.8
    6  = 115.2
    7  = 187.2
    8  = 100.8
    9  = 100.8
    10 = 144
    11 = 129.6
    12 = 100.8
    13 = 115.2
    14 = 86.4
    15 = 158.4
    16 = 86.4
    17 = 86.4
    18 = 86.4
    19 = 100.8
    20 = 129.6
    21 = 72
    22 = 115.2
    23 = 244.8
    24 = 100.8
    25 = 115.2
    26 = 115.2
    27 = 115.2
    28 = 129.6
    29 = 86.4
    30 = 115.2
    31 = 86.4
    32 = 86.4
    33 = 86.4
    34 = 144
    35 = 187.2
    36 = 115.2
    37 = 72
    38 = 28.8
    39 = 28.8
    40 = 129.6
    41 = 57.6
    42 = 86.4
}

foreach ($r in $rowHeights.Keys) {
    $ws.Rows.Item($r).RowHeight = $rowHeights[$r]
}

# ---------------------------------------------------------------------------
# Match the author's final cursor position/selection.
# ---------------------------------------------------------------------------
$ws.Range("B36").Select()
